$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the execution/test dates in column D for rows 2 and 3
$ws.Range("D2").Value = 40179
$ws.Range("D3").Value = 40179

# Update the active selection to match the recorded view state
$ws.Range("F6").Select()
